$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value2 = 30000
$ws.Range("I21").Value2 = 30000
$ws.Range("K21").Value2 = 30000
$ws.Range("M21").Value2 = -29532
$ws.Range("H23").Value2 = 30000
$ws.Range("I23").Value2 = 30000
$ws.Range("K23").Value2 = 30000
$ws.Range("M23").Value2 = -29766
$ws.Range("H58").Value2 = 1801.4445
$ws.Range("I58").Value2 = 1332.5555
$ws.Range("J58").Value2 = 2270.3333
$ws.Range("K58").Value2 = 3997.6665
$ws.Range("L58").Value2 = 6810.999899999999
$ws.Range("M58").Value2 = -3847.6665
$ws.Range("N58").Value2 = -7110.999899999999
$ws.Range("H86").Value2 = 455
$ws.Range("I86").Value2 = 457.5
$ws.Range("J86").Value2 = 450
$ws.Range("K86").Value2 = 457.5
$ws.Range("L86").Value2 = 450
$ws.Range("M86").Value2 = 665.5
$ws.Range("N86").Value2 = -2696
$ws.Range("H87").Value2 = 99999
$ws.Range("J87").Value2 = 99999
$ws.Range("L87").Value2 = 99999
$ws.Range("N87").Value2 = -102495
$ws.Range("H89").Value2 = 455
$ws.Range("I89").Value2 = 457.5
$ws.Range("J89").Value2 = 450
$ws.Range("K89").Value2 = 2287.5
$ws.Range("L89").Value2 = 2250
$ws.Range("M89").Value2 = 3328.5
$ws.Range("N89").Value2 = -13482
$ws.Range("H90").Value2 = 99999
$ws.Range("J90").Value2 = 99999
$ws.Range("L90").Value2 = 299997
$ws.Range("N90").Value2 = -312477
$ws.Range("H113").Value2 = 3721.1428
$ws.Range("I113").Value2 = 3721.1428
$ws.Range("K113").Value2 = 3721.1428
$ws.Range("M113").Value2 = -467.1428000000001
$ws.Range("H116").Value2 = 3898.5
$ws.Range("J116").Value2 = 0
$ws.Range("L116").Value2 = 0
$ws.Range("N116").ClearContents()
$ws.Range("H137").Value2 = 2614.76
$ws.Range("I137").Value2 = 1552.8572
$ws.Range("J137").Value2 = 3027.7222
$ws.Range("K137").Value2 = 4658.571599999999
$ws.Range("L137").Value2 = 9083.1666
$ws.Range("M137").Value2 = -2108.571599999999
$ws.Range("N137").Value2 = -14183.1666

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 970.9286
$ws.Range("I2").Value2 = 970.9286
$ws.Range("J2").Value2 = 0
$ws.Range("K2").Value2 = 970.9286
$ws.Range("L2").Value2 = 0
$ws.Range("M2").Value2 = -857.9286
$ws.Range("N2").ClearContents()
$ws.Range("H16").Value2 = 25498
$ws.Range("I16").Value2 = 25498
$ws.Range("K16").Value2 = 25498
$ws.Range("M16").Value2 = -25211
$ws.Range("H36").Value2 = 0
$ws.Range("I36").Value2 = 0
$ws.Range("K36").Value2 = 0
$ws.Range("M36").ClearContents()
$ws.Range("H61").Value2 = 8250
$ws.Range("I61").Value2 = 7500
$ws.Range("K61").Value2 = 7500
$ws.Range("M61").Value2 = -7288
$ws.Range("H116").Value2 = 970.9286
$ws.Range("I116").Value2 = 970.9286
$ws.Range("J116").Value2 = 0
$ws.Range("K116").Value2 = 970.9286
$ws.Range("L116").Value2 = 0
$ws.Range("M116").Value2 = 1323.0714
$ws.Range("N116").ClearContents()
$ws.Range("H136").Value2 = 8250
$ws.Range("I136").Value2 = 7500
$ws.Range("K136").Value2 = 22500
$ws.Range("M136").Value2 = -19950

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 970.9286
$ws.Range("I3").Value2 = 970.9286
$ws.Range("J3").Value2 = 0
$ws.Range("K3").Value2 = 970.9286
$ws.Range("L3").Value2 = 0
$ws.Range("M3").Value2 = -856.9286
$ws.Range("N3").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value2 = 0
$ws.Range("I17").Value2 = 0
$ws.Range("J17").Value2 = 0
$ws.Range("K17").Value2 = 0
$ws.Range("L17").Value2 = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()
$ws.Range("H36").Value2 = 100
$ws.Range("I36").Value2 = 100
$ws.Range("K36").Value2 = 100
$ws.Range("M36").Value2 = 288
$ws.Range("H40").Value2 = 100
$ws.Range("I40").Value2 = 100
$ws.Range("K40").Value2 = 100
$ws.Range("M40").Value2 = 60
$ws.Range("H50").Value2 = 40118.418
$ws.Range("I50").Value2 = 12866.6
$ws.Range("J50").Value2 = 59584
$ws.Range("K50").Value2 = 12866.6
$ws.Range("L50").Value2 = 59584
$ws.Range("M50").Value2 = -12241.6
$ws.Range("N50").Value2 = -60834
$ws.Range("H132").Value2 = 1680.2858
$ws.Range("I132").Value2 = 1680.2858
$ws.Range("K132").Value2 = 5040.857400000001
$ws.Range("M132").Value2 = -2510.857400000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value2 = 11030.5
$ws.Range("I9").Value2 = 0
$ws.Range("J9").Value2 = 11030.5
$ws.Range("K9").Value2 = 0
$ws.Range("L9").Value2 = 33091.5
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value2 = -33539.5
$ws.Range("H55").Value2 = 5966.6
$ws.Range("I55").Value2 = 1500
$ws.Range("J55").Value2 = 7083.25
$ws.Range("K55").Value2 = 4500
$ws.Range("L55").Value2 = 21249.75
$ws.Range("M55").Value2 = -4323
$ws.Range("N55").Value2 = -21603.75
$ws.Range("H116").Value2 = 2249.25
$ws.Range("I116").Value2 = 2249.25
$ws.Range("K116").Value2 = 6747.75
$ws.Range("M116").Value2 = -3305.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 1871.75
$ws.Range("J7").Value2 = 2745
$ws.Range("L7").Value2 = 2745
$ws.Range("N7").Value2 = -2969
$ws.Range("H13").Value2 = 1000
$ws.Range("I13").Value2 = 1000
$ws.Range("K13").Value2 = 1000
$ws.Range("M13").Value2 = -860
$ws.Range("H42").Value2 = 39999
$ws.Range("I42").Value2 = 0
$ws.Range("J42").Value2 = 39999
$ws.Range("K42").Value2 = 0
$ws.Range("L42").Value2 = 39999
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value2 = -41125
$ws.Range("H49").Value2 = 39999
$ws.Range("I49").Value2 = 0
$ws.Range("J49").Value2 = 39999
$ws.Range("K49").Value2 = 0
$ws.Range("L49").Value2 = 39999
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value2 = -40293
$ws.Range("H126").Value2 = 1871.75
$ws.Range("J126").Value2 = 2745
$ws.Range("L126").Value2 = 8235
$ws.Range("N126").Value2 = -13175

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value2 = 1396.2
$ws.Range("I132").Value2 = 729
$ws.Range("J132").Value2 = 2953
$ws.Range("K132").Value2 = 2187
$ws.Range("L132").Value2 = 8859
$ws.Range("M132").Value2 = 343
$ws.Range("N132").Value2 = -13919
